$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Set column A width (character width units; value chosen so the
# stored OOXML width resolves to exactly 20)
$ws.Columns.Item(1).ColumnWidth = 19.1666666666667

# Update row 2 values (B2 becomes numeric, C2/D2 updated values)
$ws.Cells.Item(2, 2).Value = 13.98
$ws.Cells.Item(2, 3).Value = 0.0003301233676692334
$ws.Cells.Item(2, 4).Value = 0.00468879702841619

# Add new row 3
$ws.Cells.Item(3, 1).Value = "23 Mar 22, 13:58PM"
$ws.Cells.Item(3, 2).Value = 13.98
$ws.Cells.Item(3, 3).Value = 0.0003301233676692334
$ws.Cells.Item(3, 4).Value = 0.00468879702841619
